$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 3323
$ws.Range("I3").Value = 3416
$ws.Range("F4").Value = 1859
$ws.Range("I4").Value = 807
$ws.Range("I5").Value = 314
$ws.Range("I6").Value = 3856
$ws.Range("F7").Value = 24048
$ws.Range("I7").Value = 11716

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("I2").Value = 33
$ws.Range("I4").Value = 18
$ws.Range("I7").Value = 133

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I3").Value = 114
$ws.Range("I6").Value = 104
$ws.Range("I7").Value = 376

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I2").Value = 109
$ws.Range("I7").Value = 459

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I2").Value = 36
$ws.Range("I7").Value = 109

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("I2").Value = 39
$ws.Range("I7").Value = 103

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I7").Value = 389
$ws.Range("I9").Value = 56
$ws.Range("I11").Value = 187
$ws.Range("I17").Value = 13
$ws.Range("I20").Value = 297
$ws.Range("I21").Value = 64
$ws.Range("I23").Value = 108
$ws.Range("I26").Value = 14
$ws.Range("I27").Value = 103
$ws.Range("I29").Value = 760
$ws.Range("I31").Value = 109
$ws.Range("I33").Value = 521
$ws.Range("I36").Value = 158
$ws.Range("I37").Value = 376
$ws.Range("I42").Value = 406
$ws.Range("I43").Value = 104
$ws.Range("I44").Value = 89
$ws.Range("I52").Value = 252
$ws.Range("I54").Value = 259
$ws.Range("I57").Value = 44
$ws.Range("I61").Value = 15
$ws.Range("F63").Value = 154
$ws.Range("I63").Value = 42
$ws.Range("I64").Value = 113
$ws.Range("I67").Value = 459
$ws.Range("I69").Value = 28
$ws.Range("I70").Value = 23
$ws.Range("I73").Value = 98
$ws.Range("I75").Value = 40
$ws.Range("I78").Value = 160
$ws.Range("I79").Value = 307
$ws.Range("I80").Value = 43
$ws.Range("I84").Value = 103
$ws.Range("I85").Value = 538
$ws.Range("I86").Value = 68
$ws.Range("I88").Value = 107
$ws.Range("I89").Value = 133
$ws.Range("I90").Value = 146
$ws.Range("I91").Value = 143
$ws.Range("F101").Value = 24048
$ws.Range("I101").Value = 11716

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I3").Value = 180
$ws.Range("I7").Value = 521

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I2").Value = 60
$ws.Range("I6").Value = 134
$ws.Range("I7").Value = 259

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 229
$ws.Range("I3").Value = 265
$ws.Range("I6").Value = 203
$ws.Range("I7").Value = 760

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("I2").Value = 32
$ws.Range("I7").Value = 89

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 137
$ws.Range("I3").Value = 217
$ws.Range("I7").Value = 538

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I2").Value = 109
$ws.Range("I3").Value = 140
$ws.Range("I7").Value = 406

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I4").Value = 24
$ws.Range("I7").Value = 160

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("I2").Value = 31
$ws.Range("I7").Value = 108

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("I4").Value = 4
$ws.Range("I7").Value = 28

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("I3").Value = 47
$ws.Range("I6").Value = 45
$ws.Range("I7").Value = 143

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("I4").Value = 2
$ws.Range("I6").Value = 51
$ws.Range("I7").Value = 64

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I2").Value = 89
$ws.Range("I3").Value = 99
$ws.Range("I6").Value = 91
$ws.Range("I7").Value = 307

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("I4").Value = 5
$ws.Range("I6").Value = 42
$ws.Range("I7").Value = 113

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I3").Value = 89
$ws.Range("I6").Value = 98
$ws.Range("I7").Value = 297

$ws = $wb.Worksheets.Item("Burnside")
$ws.Range("I3").Value = 5
$ws.Range("I7").Value = 13

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I3").Value = 48
$ws.Range("I6").Value = 49
$ws.Range("I7").Value = 158

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I2").Value = 71
$ws.Range("I7").Value = 252

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("I2").Value = 3
$ws.Range("I7").Value = 14

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I2").Value = 84
$ws.Range("I7").Value = 187

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("I3").Value = 20
$ws.Range("I7").Value = 56

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("I4").Value = 11
$ws.Range("I7").Value = 98

$ws = $wb.Worksheets.Item("O'Hare")
$ws.Range("I4").Value = 4
$ws.Range("I7").Value = 23

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("I2").Value = 25
$ws.Range("I3").Value = 38
$ws.Range("I7").Value = 107

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I3").Value = 21
$ws.Range("I7").Value = 103

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("I3").Value = 5
$ws.Range("I7").Value = 68

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("I3").Value = 13
$ws.Range("I7").Value = 40

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I3").Value = 31
$ws.Range("I7").Value = 146

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("I5").Value = 3
$ws.Range("I7").Value = 44

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("I6").Value = 62
$ws.Range("I7").Value = 104

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("I3").Value = 8
$ws.Range("I7").Value = 43

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I3").Value = 117
$ws.Range("I7").Value = 389

$ws = $wb.Worksheets.Item("Mount Greenwood")
$ws.Range("I2").Value = 6
$ws.Range("I7").Value = 15
